$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns B (Coin name) and C (Link) ---
$ws.Range("B42").Value = 'Aave'
$ws.Range("B43").Value = 'Filecoin'
$ws.Range("B44").Value = 'InjectiveProtocol'
$ws.Range("B45").Value = 'Hedera'
$ws.Range("B46").Value = 'Mantle'
$ws.Range("B47").Value = 'Stellar'
$ws.Range("B48").Value = 'VeChain'
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("B51").Value = 'WhiteBITCoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("C43").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C44").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("C45").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("C47").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C51").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'

# --- Column D (Price) : force text format to avoid numeric coercion ---
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"
$ws.Range("D2").Value = '63.138.04'
$ws.Range("D3").Value = '2.563.08'
$ws.Range("D5").Value = '583.84'
$ws.Range("D6").Value = '147.77'
$ws.Range("D10").Value = '5.59'
$ws.Range("D13").Value = '27.58'
$ws.Range("D14").Value = '3.021.69'
$ws.Range("D15").Value = '63.051.91'
$ws.Range("D17").Value = '2.558.31'
$ws.Range("D18").Value = '11.38'
$ws.Range("D19").Value = '341.02'
$ws.Range("D20").Value = '4.38'
$ws.Range("D21").Value = '6.83'
$ws.Range("D24").Value = '2.681.14'
$ws.Range("D28").Value = '8.45'
$ws.Range("D30").Value = '7.87'
$ws.Range("D31").Value = '1.97'
$ws.Range("D32").Value = '0.0₃0825'
$ws.Range("D33").Value = '176.63'
$ws.Range("D35").Value = '432.69'
$ws.Range("D37").Value = '19.23'
$ws.Range("D38").Value = '4.49'
$ws.Range("D41").Value = '0.999'
$ws.Range("D42").Value = '151.93'
$ws.Range("D43").Value = '3.81'
$ws.Range("D44").Value = '21.06'
$ws.Range("D45").Value = '0.0551'
$ws.Range("D46").Value = '0.606'
$ws.Range("D47").Value = '0.0974'
$ws.Range("D48").Value = '0.0242'
$ws.Range("D49").Value = '18.35'
$ws.Range("D50").Value = '1.72'
$ws.Range("D51").Value = '11.37'
$colD.Style = "Normal"

# --- Column E (Volume/1h %) : force text format to avoid numeric coercion ---
$colE = $ws.Range("E2:E51")
$colE.NumberFormat = "@"
$ws.Range("E2").Value = '  +0.52%  '
$ws.Range("E3").Value = '  +1.24%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("E5").Value = '  +2.95%  '
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("E8").Value = '  +0.71%  '
$ws.Range("E9").Value = '  +1.94%  '
$ws.Range("E10").Value = '  -0.83%  '
$ws.Range("E11").Value = '  +0.30%  '
$ws.Range("E12").Value = '  +0.34%  '
$ws.Range("E13").Value = '  -0.52%  '
$ws.Range("E14").Value = '  +1.20%  '
$ws.Range("E15").Value = '  +0.34%  '
$ws.Range("E16").Value = '  +2.44%  '
$ws.Range("E17").Value = '  +0.49%  '
$ws.Range("E18").Value = '  -1.56%  '
$ws.Range("E19").Value = '  +1.45%  '
$ws.Range("E20").Value = '  +2.13%  '
$ws.Range("E21").Value = '  +1.01%  '
$ws.Range("E22").Value = '  +0.05%  '
$ws.Range("E23").Value = '  +0.53%  '
$ws.Range("E24").Value = '  +0.92%  '
$ws.Range("E25").Value = '  +3.17%  '
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").Value = '  +1.34%  '
$ws.Range("E29").Value = '  -2.58%  '
$ws.Range("E30").Value = '  +8.88%  '
$ws.Range("E31").Value = '  +6.71%  '
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("E33").Value = '  -0.56%  '
$ws.Range("E34").Value = '  +1.42%  '
$ws.Range("E35").Value = '  +5.25%  '
$ws.Range("E36").Value = '  +1.24%  '
$ws.Range("E37").Value = '  +2.05%  '
$ws.Range("E38").Value = '  +1.93%  '
$ws.Range("E39").Value = '  +0.03%  '
$ws.Range("E40").Value = '  +0.88%  '
$ws.Range("E41").Value = '  -0.23%  '
$ws.Range("E42").Value = '  -0.30%  '
$ws.Range("E43").Value = '  +1.86%  '
$ws.Range("E44").Value = '  +1.91%  '
$ws.Range("E45").Value = '  +5.96%  '
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("E48").Value = '  +2.46%  '
$ws.Range("E49").Value = '  -0.38%  '
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  -0.32%  '
$colE.Style = "Normal"
